# Auto-generated edit script: updates Moogle Profits market-price derived
# columns (H:N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets to match
# refreshed market-board data from the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 902.1667
$ws.Range("I2").Value = 169.44444
$ws.Range("J2").Value = 3100.3333
$ws.Range("K2").Value = 169.44444
$ws.Range("L2").Value = 3100.3333
$ws.Range("M2").Value = -56.44443999999999
$ws.Range("N2").Value = -3326.3333
$ws.Range("H4").Value = 2390
$ws.Range("I4").Value = 577.5
$ws.Range("J4").Value = 3598.3333
$ws.Range("K4").Value = 577.5
$ws.Range("L4").Value = 3598.3333
$ws.Range("M4").Value = -463.5
$ws.Range("N4").Value = -3826.3333
$ws.Range("H19").Value = 1340.2693
$ws.Range("I19").Value = 998.7857
$ws.Range("J19").Value = 1738.6666
$ws.Range("K19").Value = 998.7857
$ws.Range("L19").Value = 1738.6666
$ws.Range("M19").Value = -823.7857
$ws.Range("N19").Value = -2088.6666
$ws.Range("H32").Value = 4582.5
$ws.Range("I32").Value = 3999.3333
$ws.Range("J32").Value = 4776.8887
$ws.Range("K32").Value = 3999.3333
$ws.Range("L32").Value = 4776.8887
$ws.Range("M32").Value = -3673.3333
$ws.Range("N32").Value = -5428.8887
$ws.Range("H40").Value = 2421.0908
$ws.Range("I40").Value = 1300.8
$ws.Range("K40").Value = 1300.8
$ws.Range("M40").Value = -1125.8
$ws.Range("H43").Value = 2371.3809
$ws.Range("I43").Value = 1666.3636
$ws.Range("J43").Value = 3146.9
$ws.Range("K43").Value = 1666.3636
$ws.Range("L43").Value = 3146.9
$ws.Range("M43").Value = -1597.3636
$ws.Range("N43").Value = -3284.9
$ws.Range("H69").Value = 16438.834
$ws.Range("I69").Value = 3976.25
$ws.Range("K69").Value = 11928.75
$ws.Range("M69").Value = -11054.75
$ws.Range("H72").Value = 16438.834
$ws.Range("I72").Value = 3976.25
$ws.Range("K72").Value = 35786.25
$ws.Range("M72").Value = -31418.25
$ws.Range("H132").Value = 1822.9482
$ws.Range("I132").Value = 1637.1273
$ws.Range("K132").Value = 4911.3819
$ws.Range("M132").Value = -2381.3819
$ws.Range("H135").Value = 768.3333
$ws.Range("I135").Value = 768.3333
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 6914.9997
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -4379.9997
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 1989.275
$ws.Range("I137").Value = 2046.1714
$ws.Range("J137").Value = 1591
$ws.Range("K137").Value = 6138.5142
$ws.Range("L137").Value = 4773
$ws.Range("M137").Value = -3588.5142
$ws.Range("N137").Value = -9873
$ws.Range("H138").Value = 5059.8384
$ws.Range("I138").Value = 4054.617
$ws.Range("J138").Value = 7309.619
$ws.Range("K138").Value = 12163.851
$ws.Range("L138").Value = 21928.857
$ws.Range("M138").Value = -7023.851000000001
$ws.Range("N138").Value = -32208.857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1320.1
$ws.Range("I102").Value = 1320.1
$ws.Range("K102").Value = 1320.1
$ws.Range("M102").Value = 301.9000000000001
$ws.Range("H110").Value = 1433.7778
$ws.Range("I110").Value = 1464.3636
$ws.Range("J110").Value = 1299.2
$ws.Range("K110").Value = 1464.3636
$ws.Range("L110").Value = 1299.2
$ws.Range("M110").Value = 580.6364000000001
$ws.Range("N110").Value = -5389.2
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 2705.3518
$ws.Range("I132").Value = 1708.6136
$ws.Range("K132").Value = 5125.8408
$ws.Range("M132").Value = -2595.8408

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5073.7144
$ws.Range("I105").Value = 5123.75
$ws.Range("K105").Value = 5123.75
$ws.Range("M105").Value = -3376.75
$ws.Range("H134").Value = 2424.0476
$ws.Range("I134").Value = 1242.4474
$ws.Range("J134").Value = 13649.25
$ws.Range("K134").Value = 3727.3422
$ws.Range("L134").Value = 40947.75
$ws.Range("M134").Value = -1192.3422
$ws.Range("N134").Value = -46017.75
$ws.Range("H139").Value = 83999.60000000001
$ws.Range("J139").Value = 120000
$ws.Range("L139").Value = 120000
$ws.Range("N139").Value = -130280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 980.8333
$ws.Range("I22").Value = 524.3333
$ws.Range("K22").Value = 524.3333
$ws.Range("M22").Value = -174.3333
$ws.Range("H68").Value = 79999
$ws.Range("J68").Value = 79999
$ws.Range("L68").Value = 79999
$ws.Range("N68").Value = -81497
$ws.Range("H71").Value = 79999
$ws.Range("J71").Value = 79999
$ws.Range("L71").Value = 239997
$ws.Range("N71").Value = -247485
$ws.Range("H99").Value = 1931.1818
$ws.Range("I99").Value = 1952.0526
$ws.Range("J99").Value = 1799
$ws.Range("K99").Value = 1952.0526
$ws.Range("L99").Value = 1799
$ws.Range("M99").Value = -454.0526
$ws.Range("N99").Value = -4795
$ws.Range("H107").Value = 1366.0435
$ws.Range("J107").Value = 1918.2
$ws.Range("L107").Value = 1918.2
$ws.Range("N107").Value = -5758.2
$ws.Range("H126").Value = 1931.1818
$ws.Range("I126").Value = 1952.0526
$ws.Range("J126").Value = 1799
$ws.Range("K126").Value = 5856.1578
$ws.Range("L126").Value = 5397
$ws.Range("M126").Value = -3386.1578
$ws.Range("N126").Value = -10337
$ws.Range("H132").Value = 2017.262
$ws.Range("I132").Value = 1351.9487
$ws.Range("J132").Value = 10666.333
$ws.Range("K132").Value = 4055.8461
$ws.Range("L132").Value = 31998.999
$ws.Range("M132").Value = -1525.8461
$ws.Range("N132").Value = -37058.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1350.1538
$ws.Range("I5").Value = 1591.4286
$ws.Range("J5").Value = 1068.6666
$ws.Range("K5").Value = 4774.2858
$ws.Range("L5").Value = 3205.9998
$ws.Range("M5").Value = -4662.2858
$ws.Range("N5").Value = -3429.9998
$ws.Range("H46").Value = 1048.0834
$ws.Range("I46").Value = 155.4
$ws.Range("J46").Value = 1685.7142
$ws.Range("K46").Value = 466.2
$ws.Range("L46").Value = 5057.142599999999
$ws.Range("M46").Value = -375.2
$ws.Range("N46").Value = -5239.142599999999
$ws.Range("H70").Value = 13666.833
$ws.Range("I70").Value = 1000.5
$ws.Range("K70").Value = 3001.5
$ws.Range("M70").Value = -2686.5
$ws.Range("H73").Value = 13666.833
$ws.Range("I73").Value = 1000.5
$ws.Range("K73").Value = 3001.5
$ws.Range("M73").Value = -1909.5
$ws.Range("H86").Value = 10092.833
$ws.Range("I86").Value = 25901
$ws.Range("J86").Value = 2188.75
$ws.Range("K86").Value = 77703
$ws.Range("L86").Value = 6566.25
$ws.Range("M86").Value = -76517
$ws.Range("N86").Value = -8938.25
$ws.Range("H89").Value = 10092.833
$ws.Range("I89").Value = 25901
$ws.Range("J89").Value = 2188.75
$ws.Range("K89").Value = 233109
$ws.Range("L89").Value = 19698.75
$ws.Range("M89").Value = -227181
$ws.Range("N89").Value = -31554.75
$ws.Range("H124").Value = 4029.5
$ws.Range("I124").Value = 4029.5
$ws.Range("K124").Value = 12088.5
$ws.Range("M124").Value = -7178.5
$ws.Range("H135").Value = 1350.1538
$ws.Range("I135").Value = 1591.4286
$ws.Range("J135").Value = 1068.6666
$ws.Range("K135").Value = 14322.8574
$ws.Range("L135").Value = 9617.999400000001
$ws.Range("M135").Value = -11787.8574
$ws.Range("N135").Value = -14687.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5127.5
$ws.Range("I80").Value = 3921.4119
$ws.Range("J80").Value = 6991.4546
$ws.Range("K80").Value = 3921.4119
$ws.Range("L80").Value = 6991.4546
$ws.Range("M80").Value = -2923.4119
$ws.Range("N80").Value = -8987.454600000001
$ws.Range("H83").Value = 5127.5
$ws.Range("I83").Value = 3921.4119
$ws.Range("J83").Value = 6991.4546
$ws.Range("K83").Value = 19607.0595
$ws.Range("L83").Value = 34957.273
$ws.Range("M83").Value = -14615.0595
$ws.Range("N83").Value = -44941.273
$ws.Range("H102").Value = 2692.0303
$ws.Range("I102").Value = 1398.375
$ws.Range("K102").Value = 1398.375
$ws.Range("M102").Value = 223.625
$ws.Range("H132").Value = 2111.9756
$ws.Range("I132").Value = 1539.5151
$ws.Range("K132").Value = 4618.5453
$ws.Range("M132").Value = -2088.5453

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2487.5312
$ws.Range("I16").Value = 381.5238
$ws.Range("J16").Value = 6508.091
$ws.Range("K16").Value = 381.5238
$ws.Range("L16").Value = 6508.091
$ws.Range("M16").Value = -211.5238
$ws.Range("N16").Value = -6848.091
$ws.Range("H21").Value = 14333.333
$ws.Range("J21").Value = 14333.333
$ws.Range("L21").Value = 14333.333
$ws.Range("N21").Value = -14681.333
$ws.Range("H22").Value = 2387.1428
$ws.Range("I22").Value = 2450
$ws.Range("K22").Value = 2450
$ws.Range("M22").Value = -2155
$ws.Range("H27").Value = 2387.1428
$ws.Range("I27").Value = 2450
$ws.Range("K27").Value = 2450
$ws.Range("M27").Value = -2343
$ws.Range("H46").Value = 2096
$ws.Range("I46").Value = 890.6667
$ws.Range("J46").Value = 3452
$ws.Range("K46").Value = 890.6667
$ws.Range("L46").Value = 3452
$ws.Range("M46").Value = -702.6667
$ws.Range("N46").Value = -3828
$ws.Range("H55").Value = 705.48
$ws.Range("J55").Value = 1217.6154
$ws.Range("L55").Value = 1217.6154
$ws.Range("N55").Value = -1563.6154
$ws.Range("H82").Value = 1378.6666
$ws.Range("I82").Value = 907.44446
$ws.Range("K82").Value = 907.44446
$ws.Range("M82").Value = -546.44446
$ws.Range("H85").Value = 1378.6666
$ws.Range("I85").Value = 907.44446
$ws.Range("K85").Value = 907.44446
$ws.Range("M85").Value = 340.55554
$ws.Range("H132").Value = 4933.273
$ws.Range("I132").Value = 2727.6667
$ws.Range("K132").Value = 8183.000100000001
$ws.Range("M132").Value = -5653.000100000001
